$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Handle interference without PUN": process times for therblig rows
# A (row 55) and DA (row 56) increase from 120 to 138 across LH/RH/BOT.
$ws.Range("B55:D55").Value = 138
$ws.Range("B56:D56").Value = 138

# Update the saved view state to match where the author left off editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("B55:D56").Select()
